# Revert "feat: update config" — restore the original Guid values in
# column C (rows 19-31) that a previous commit had overwritten with
# placeholder date-serial-looking numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = 136198
$ws.Range("C20").Value = 203458
$ws.Range("C21").Value = 203458
$ws.Range("C22").Value = 134424
$ws.Range("C23").Value = 204675
$ws.Range("C24").Value = 169179
$ws.Range("C25").Value = 136198
$ws.Range("C26").Value = 169139
$ws.Range("C27").Value = 169180
$ws.Range("C28").Value = 199679
$ws.Range("C29").Value = 206005
$ws.Range("C30").Value = 204923
$ws.Range("C31").Value = 169189

# Restore the cursor/selection position recorded in the target file.
$ws.Range("E36").Select()
